$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10-16 block (MEC-2A-Ajustagem shifting from column E to F)
$ws.Range("F10").Value = "[-, -, -, 'MEC-2A-Ajustagem']"

$ws.Range("E11").Value = "-"
$ws.Range("E12").Value = "-"

$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "['MEC-2A-Ajustagem', -, -, -]"

$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "['MEC-2A-Ajustagem', -, -, -]"

$ws.Range("F16").Value = "['MEC-2A-Ajustagem', -, -, -]"

# Row 18
$ws.Range("B18").Value = "[-, 'MEC-1NA-Tornearia', 'MEC-1NB-Tornearia', 'MEC-1NA-Metrologia 1']"
$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "[-, -, 'ELM-1NA-Processos de Usinagem 1', -]"
$ws.Range("F18").Value = "-"

# Row 19
$ws.Range("B19").Value = "[-, 'MEC-1NA-Tornearia', 'MEC-1NB-Tornearia', 'MEC-1NA-Metrologia 1']"
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = "[-, -, 'ELM-1NA-Processos de Usinagem 1', -]"
$ws.Range("F19").Value = "-"

# Row 20
$ws.Range("B20").Value = "[-, -, 'MEC-1NB-Tornearia', -]"
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "[-, 'MEC-1NA-Tornearia', -, 'MEC-1NA-Metrologia 1']"
$ws.Range("E20").Value = "[-, -, 'ELM-1NA-Processos de Usinagem 1', -]"
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("B21").Value = "[-, -, 'MEC-1NB-Tornearia', -]"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "[-, 'MEC-1NA-Tornearia', -, 'MEC-1NA-Metrologia 1']"
$ws.Range("E21").Value = "[-, -, 'ELM-1NA-Processos de Usinagem 1', -]"
$ws.Range("F21").Value = "-"

$wb.Save()
